# Applies the crypto price/volume refresh described by the commit diff.
# Price-column values that are valid numeric literals (e.g. "268.18") are
# written with a leading apostrophe so Excel keeps them as literal TEXT
# (matching the inline-string cells in the workbook) instead of silently
# converting them to floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.704.95'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '2.284.26'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '''114.92'
$ws.Range("E5").Value = '  +10.66%  '
$ws.Range("D6").Value = '''268.18'
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("D7").Value = '''0.631'
$ws.Range("E7").Value = '  +2.29%  '
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").Value = '''0.621'
$ws.Range("E9").Value = '  +2.35%  '
$ws.Range("D10").Value = '''49.07'
$ws.Range("E10").Value = '  +6.79%  '
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("E12").Value = '  +13.07%  '
$ws.Range("E13").Value = '  +0.56%  '
$ws.Range("D14").Value = '''15.83'
$ws.Range("E14").Value = '  +1.41%  '
$ws.Range("D15").Value = '2.629.69'
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").Value = '''0.884'
$ws.Range("E16").Value = '  +3.16%  '
$ws.Range("D17").Value = '2.287.69'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").Value = '43.592.37'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = '''0.0000109'
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("D20").Value = '''7.00'
$ws.Range("E20").Value = '  +11.87%  '
$ws.Range("D21").Value = '''72.32'
$ws.Range("E22").Value = '  -2.69%  '
$ws.Range("D23").Value = '''9.92'
$ws.Range("E23").Value = '  +8.47%  '
$ws.Range("D24").Value = '''233.35'
$ws.Range("D25").Value = '''2.89'
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = '''11.64'
$ws.Range("E27").Value = '  +3.87%  '
$ws.Range("D28").Value = '''41.99'
$ws.Range("E28").Value = '  +4.43%  '
$ws.Range("E29").Value = '  -1.68%  '
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").Value = '''173.24'
$ws.Range("E31").Value = '  -2.19%  '
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("D33").Value = '''0.0925'
$ws.Range("E33").Value = '  +2.82%  '
$ws.Range("D34").Value = '''5.73'
$ws.Range("E34").Value = '  +4.94%  '
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("D36").Value = '''4.67'
$ws.Range("E36").Value = '  -4.97%  '
$ws.Range("D37").Value = '''0.0355'
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("D38").Value = '''0.106'
$ws.Range("E38").Value = '  -2.83%  '
$ws.Range("D39").Value = '''3.76'
$ws.Range("E39").Value = '  +6.03%  '
$ws.Range("D40").Value = '''14.95'
$ws.Range("E40").Value = '  +21.72%  '
$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").Value = '''75.07'
$ws.Range("E41").Value = '  +14.75%  '
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").Value = '''2.44'
$ws.Range("E42").Value = '  +4.48%  '
$ws.Range("D43").Value = '''0.239'
$ws.Range("E43").Value = '  +0.41%  '
$ws.Range("D44").Value = '''6.32'
$ws.Range("E44").Value = '  +20.48%  '
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").Value = '''1.38'
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").Value = '''1.26'
$ws.Range("E48").Value = '  +3.13%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '''102.69'
$ws.Range("E49").Value = '  +3.82%  '
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").Value = '''0.456'
$ws.Range("E51").Value = '  +1.73%  '
